# Generate Report for Handback
# Update the Correspond Handoff/Handback Datetime values on the
# per-locale sheets (zh-cn, de-de). Rows 2 and 3 on each sheet shared the
# same handoff/handback timestamps, so both rows are refreshed together.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 02:14:55"
$wsZhCn.Range("H2").Value = "2016-03-14 02:15:10"
$wsZhCn.Range("E3").Value = "2016-03-14 02:14:55"
$wsZhCn.Range("H3").Value = "2016-03-14 02:15:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 02:14:58"
$wsDeDe.Range("H2").Value = "2016-03-14 02:15:15"
$wsDeDe.Range("E3").Value = "2016-03-14 02:14:58"
$wsDeDe.Range("H3").Value = "2016-03-14 02:15:15"
